# Auto-generated edit script: apply scheduled market-data refresh to Sheets
# Updates currentAveragePrice / Leve profit columns (H-N) per the commit diff.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 453
$ws.Range("I6").Value = 316.25
$ws.Range("K6").Value = 948.75
$ws.Range("M6").Value = -836.75
$ws.Range("H32").Value = 388.125
$ws.Range("I32").Value = 416.6
$ws.Range("J32").Value = 340.66666
$ws.Range("K32").Value = 416.6
$ws.Range("L32").Value = 340.66666
$ws.Range("M32").Value = -90.60000000000002
$ws.Range("N32").Value = -992.66666
$ws.Range("H98").Value = 1986.4286
$ws.Range("I98").Value = 1394.6471
$ws.Range("K98").Value = 1394.6471
$ws.Range("M98").Value = 103.3529000000001
$ws.Range("H106").Value = 4373.294
$ws.Range("I106").Value = 3326.6667
$ws.Range("J106").Value = 4944.1816
$ws.Range("K106").Value = 3326.6667
$ws.Range("L106").Value = 4944.1816
$ws.Range("M106").Value = -2695.6667
$ws.Range("N106").Value = -6206.1816
$ws.Range("H122").Value = 1986.4286
$ws.Range("I122").Value = 1394.6471
$ws.Range("K122").Value = 4183.9413
$ws.Range("M122").Value = -1733.9413
$ws.Range("H137").Value = 587160.8
$ws.Range("I137").Value = 2572.2173
$ws.Range("J137").Value = 982617.9
$ws.Range("K137").Value = 7716.651899999999
$ws.Range("L137").Value = 2947853.7
$ws.Range("M137").Value = -5166.651899999999
$ws.Range("N137").Value = -2952953.7
$ws.Range("H138").Value = 3515.4807
$ws.Range("I138").Value = 2941.7778
$ws.Range("J138").Value = 3635.558
$ws.Range("K138").Value = 8825.3334
$ws.Range("L138").Value = 10906.674
$ws.Range("M138").Value = -3685.3334
$ws.Range("N138").Value = -21186.674

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 12607.75
$ws.Range("J43").Value = 12607.75
$ws.Range("L43").Value = 12607.75
$ws.Range("N43").Value = -13233.75
$ws.Range("H97").Value = 903.7917
$ws.Range("I97").Value = 858.6818
$ws.Range("J97").Value = 1400
$ws.Range("K97").Value = 858.6818
$ws.Range("L97").Value = 1400
$ws.Range("M97").Value = -362.6818
$ws.Range("N97").Value = -2392
$ws.Range("H112").Value = 39382.09
$ws.Range("J112").Value = 39382.09
$ws.Range("L112").Value = 39382.09
$ws.Range("N112").Value = -42336.09
$ws.Range("H138").Value = 52750
$ws.Range("J138").Value = 52750
$ws.Range("L138").Value = 52750
$ws.Range("N138").Value = -63030

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 70684
$ws.Range("J108").Value = 70684
$ws.Range("L108").Value = 70684
$ws.Range("N108").Value = -78364

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2436.3845
$ws.Range("I132").Value = 2167.25
$ws.Range("J132").Value = 5666
$ws.Range("K132").Value = 6501.75
$ws.Range("L132").Value = 16998
$ws.Range("M132").Value = -3971.75
$ws.Range("N132").Value = -22058
$ws.Range("H134").Value = 1956.3478
$ws.Range("I134").Value = 1535.7142
$ws.Range("J134").Value = 3294.7273
$ws.Range("K134").Value = 4607.142599999999
$ws.Range("L134").Value = 9884.1819
$ws.Range("M134").Value = -2072.142599999999
$ws.Range("N134").Value = -14954.1819

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 42
$ws.Range("I2").Value = 35
$ws.Range("J2").Value = 42.56
$ws.Range("K2").Value = 210
$ws.Range("L2").Value = 255.36
$ws.Range("M2").Value = -97
$ws.Range("N2").Value = -481.36
$ws.Range("H9").Value = 42578.652
$ws.Range("J9").Value = 46619
$ws.Range("L9").Value = 139857
$ws.Range("N9").Value = -140305
$ws.Range("H19").Value = 500
$ws.Range("J19").Value = 500
$ws.Range("L19").Value = 1500
$ws.Range("N19").Value = -1848
$ws.Range("H22").Value = 71430424
$ws.Range("I22").Value = 142858000
$ws.Range("J22").Value = 2842.8572
$ws.Range("K22").Value = 428574000
$ws.Range("L22").Value = 8528.571599999999
$ws.Range("M22").Value = -428573831
$ws.Range("N22").Value = -8866.571599999999
$ws.Range("H27").Value = 71430424
$ws.Range("I27").Value = 142858000
$ws.Range("J27").Value = 2842.8572
$ws.Range("K27").Value = 428574000
$ws.Range("L27").Value = 8528.571599999999
$ws.Range("M27").Value = -428573898
$ws.Range("N27").Value = -8732.571599999999
$ws.Range("H33").Value = 131
$ws.Range("J33").Value = 168
$ws.Range("L33").Value = 1008
$ws.Range("N33").Value = -1574
$ws.Range("H68").Value = 2991.2932
$ws.Range("J68").Value = 5291.3335
$ws.Range("L68").Value = 15874.0005
$ws.Range("N68").Value = -17496.0005
$ws.Range("H71").Value = 2991.2932
$ws.Range("J71").Value = 5291.3335
$ws.Range("L71").Value = 47622.0015
$ws.Range("N71").Value = -55734.0015
$ws.Range("H129").Value = 1742.7179
$ws.Range("I129").Value = 1875
$ws.Range("J129").Value = 1650.6957
$ws.Range("K129").Value = 5625
$ws.Range("L129").Value = 4952.0871
$ws.Range("M129").Value = -625
$ws.Range("N129").Value = -14952.0871
$ws.Range("H131").Value = 1071.65
$ws.Range("I131").Value = 1441.3572
$ws.Range("J131").Value = 993.2273
$ws.Range("K131").Value = 4324.071599999999
$ws.Range("L131").Value = 2979.6819
$ws.Range("M131").Value = 715.9284000000007
$ws.Range("N131").Value = -13059.6819
$ws.Range("H134").Value = 3371.5386
$ws.Range("I134").Value = 3314.7778
$ws.Range("J134").Value = 3499.25
$ws.Range("K134").Value = 9944.3334
$ws.Range("L134").Value = 10497.75
$ws.Range("M134").Value = -4874.3334
$ws.Range("N134").Value = -20637.75
$ws.Range("H139").Value = 1905960.4
$ws.Range("I139").Value = 3355195.2
$ws.Range("J139").Value = 3839.875
$ws.Range("K139").Value = 10065585.6
$ws.Range("L139").Value = 11519.625
$ws.Range("M139").Value = -10060445.6
$ws.Range("N139").Value = -21799.625
$ws.Range("H140").Value = 1408.6976
$ws.Range("I140").Value = 991.65717
$ws.Range("J140").Value = 3233.25
$ws.Range("K140").Value = 2974.97151
$ws.Range("L140").Value = 9699.75
$ws.Range("M140").Value = 2205.02849
$ws.Range("N140").Value = -20059.75

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 369166660
$ws.Range("I11").Value = 369166660
$ws.Range("K11").Value = 369166660
$ws.Range("M11").Value = -369166521
$ws.Range("H46").Value = 28950
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 28950
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 28950
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -29262
$ws.Range("H48").Value = 15000
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3098.8057
$ws.Range("I40").Value = 2986.6365
$ws.Range("J40").Value = 3275.0715
$ws.Range("K40").Value = 2986.6365
$ws.Range("L40").Value = 3275.0715
$ws.Range("M40").Value = -2850.6365
$ws.Range("N40").Value = -3547.0715
$ws.Range("H110").Value = 79800
$ws.Range("J110").Value = 79800
$ws.Range("L110").Value = 79800
$ws.Range("N110").Value = -87980

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").Value = 0
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("N49").ClearContents()
